# Applies updated "number of dwellings" values for AB scenario (years 1901-2050,
# rows 303-452) to reflect the finished calibration run and the start of the
# renovation-change scenario run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B (rows 303-452, corresponding to years 1901-2050),
# taken from the updated baseline model output.
$newValues = @(
    1562.231466726676,
    3149.240758891025,
    4761.455215897783,
    6399.307790093049,
    8063.236924101251,
    9753.68641156211,
    11471.10524059318,
    13215.9474193051,
    14988.67178270165,
    16789.74178030778,
    18619.6252438775,
    20478.79413454716,
    22367.72426881849,
    24286.89502277552,
    26236.78901396181,
    28217.89176037541,
    30230.69131606622,
    32275.67788285805,
    34353.34339775828,
    36464.18109565748,
    38608.68504697403,
    40787.34966994636,
    43000.66921733767,
    45249.13723737408,
    47533.24600880643,
    49853.48595005651,
    52210.34500248212,
    54604.30798787673,
    57035.85594040353,
    59505.46541325291,
    62013.60776040423,
    64560.74839397207,
    67147.3460177168,
    69773.8518374049,
    72440.70874881349,
    75148.35050428357,
    77897.20085884085,
    80687.6726970223,
    83520.16714166098,
    86395.07264600335,
    89312.76407065538,
    92273.60174697166,
    95277.93052862627,
    98326.07883321807,
    101418.3576758894,
    104555.0596970436,
    107736.4581863681,
    110962.806105471,
    114234.33511155,
    117551.2545846056,
    120913.7506608076,
    124321.9852747057,
    127776.0952130563,
    131276.1911831023,
    134822.3568982082,
    138905.4527780479,
    143041.2135274808,
    147229.6519363095,
    151470.7470686474,
    155764.4434059678,
    160110.6500308372,
    164509.2398546874,
    168960.0488929176,
    173462.8755905561,
    178017.4802016286,
    182623.5842252702,
    187280.8699015069,
    191988.9797694895,
    196747.5162908118,
    201556.0415403774,
    205654.4506582586,
    209793.6205756541,
    213973.0403188611,
    218192.1626352042,
    222450.4039863175,
    226747.1446064895,
    231081.7286268641,
    235453.4642660647,
    239861.6240875848,
    244305.4453240507,
    247542.0213823364,
    250801.4648262378,
    254083.0435117222,
    257385.9968319098,
    260709.5360696641,
    264052.8447942914,
    267415.0793014182,
    270795.3690949908,
    274192.8174102352,
    277606.5017763217,
    283307.9551854341,
    289039.9962735071,
    294801.3477246399,
    300590.7073699429,
    306406.7497169889,
    312248.1275298574,
    318113.4734563555,
    324001.4016989213,
    329910.5097256748,
    335839.3800180302,
    351859.6762067491,
    367952.9128960276,
    384116.5541667915,
    400348.0852220349,
    416645.0180621215,
    433004.8971705797,
    449425.3051955232,
    465903.8686113079,
    482438.2633446286,
    499026.2203486587,
    514847.7654637888,
    530714.8940952768,
    546625.5831641678,
    562577.8809572249,
    578569.9119169221,
    594599.8812759762,
    610666.0795155464,
    626766.8866251045,
    642900.7761407595,
    659066.3189374893,
    675262.1867493393,
    691487.1553901108,
    707740.1076455261,
    724020.035806153,
    740326.0438087285,
    756657.3489517822,
    773013.2831497326,
    789393.2936879841,
    805796.9434399192,
    822223.9105052016,
    838673.9872275017,
    855147.0785486251,
    871643.1996552453,
    888162.4728739177,
    904705.1237700412,
    921271.476406802,
    937861.9477211658,
    954477.0409755436,
    971117.3382461244,
    987783.4919119548,
    1004476.215112791,
    1021196.27114867,
    1037944.461799973,
    1054721.614553673,
    1071528.568729396,
    1088366.160507993,
    1105235.206875488,
    1122136.488506451,
    1139070.731623216,
    1156038.588880617
)

$startRow = 303
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $newValues[$i]
}
